# Add the new "Napredovanje_izpad" worksheet as the last tab in the workbook
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Napredovanje_izpad"

# Match the outline summary settings used by the other sheets in this workbook
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1

# Header row
$ws.Range("B1").Value = "Sezona"
$ws.Range("C1").Value = "Izpadle ekipe"
$ws.Range("D1").Value = "Napredovane ekipe"
$headerRange = $ws.Range("B1:D1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Data rows: season index, season label, relegated teams, promoted teams
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "1993-1994"
$ws.Cells.Item(2, 3).Value = "{'Lerida', 'Osasuna', 'Rayo Vallecano'}"
$ws.Cells.Item(2, 4).Value = "{'Espanyol', 'Real Betis', 'Compostela'}"
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "1994-1995"
$ws.Cells.Item(3, 3).Value = "{'CD Logroñés'}"
$ws.Cells.Item(3, 4).Value = "{'Mérida UD', 'UD Salamanca', 'Rayo Vallecano'}"
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "1995-1996"
$ws.Cells.Item(4, 3).Value = "{'Albacete', 'Mérida UD', 'UD Salamanca'}"
$ws.Cells.Item(4, 4).Value = "{'CF Extremadura', 'Hércules', 'CD Logroñés'}"
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "1996-1997"
$ws.Cells.Item(5, 3).Value = "{'CD Logroñés', 'Sevilla', 'Hércules', 'Rayo Vallecano', 'CF Extremadura'}"
$ws.Cells.Item(5, 4).Value = "{'Mérida UD', 'UD Salamanca', 'Mallorca'}"
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "1997-1998"
$ws.Cells.Item(6, 3).Value = "{'Mérida UD', 'Sporting Gijon', 'Compostela'}"
$ws.Cells.Item(6, 4).Value = "{'Alaves', 'Villarreal', 'CF Extremadura'}"
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = "1998-1999"
$ws.Cells.Item(7, 3).Value = "{'Villarreal', 'Tenerife', 'UD Salamanca', 'CF Extremadura'}"
$ws.Cells.Item(7, 4).Value = "{'Numancia', 'Sevilla', 'Malaga', 'Rayo Vallecano'}"
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "1999-2000"
$ws.Cells.Item(8, 3).Value = "{'Atletico Madrid', 'Sevilla', 'Real Betis'}"
$ws.Cells.Item(8, 4).Value = "{'Las Palmas', 'Villarreal', 'Osasuna'}"
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "2000-2001"
$ws.Cells.Item(9, 3).Value = "{'Numancia', 'Real Oviedo', 'Racing Santander'}"
$ws.Cells.Item(9, 4).Value = "{'Sevilla', 'Tenerife', 'Real Betis'}"
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "2001-2002"
$ws.Cells.Item(10, 3).Value = "{'Las Palmas', 'Real Zaragoza', 'Tenerife'}"
$ws.Cells.Item(10, 4).Value = "{'Atletico Madrid', 'Racing Santander', 'Recreativo Huelva'}"
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "2002-2003"
$ws.Cells.Item(11, 3).Value = "{'Recreativo Huelva', 'Alaves', 'Rayo Vallecano'}"
$ws.Cells.Item(11, 4).Value = "{'Real Murcia', 'Albacete', 'Real Zaragoza'}"
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "2003-2004"
$ws.Cells.Item(12, 3).Value = "{'Real Murcia', 'Real Valladolid', 'Celta Vigo'}"
$ws.Cells.Item(12, 4).Value = "{'Numancia', 'Getafe', 'Levante'}"
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "2004-2005"
$ws.Cells.Item(13, 3).Value = "{'Albacete', 'Numancia', 'Levante'}"
$ws.Cells.Item(13, 4).Value = "{'Alaves', 'Cadiz', 'Celta Vigo'}"
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "2005-2006"
$ws.Cells.Item(14, 3).Value = "{'Alaves', 'Cadiz', 'Malaga'}"
$ws.Cells.Item(14, 4).Value = "{'Recreativo Huelva', 'Gimnàstic Tarr…', 'Levante'}"
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "2006-2007"
$ws.Cells.Item(15, 3).Value = "{'Gimnàstic Tarr…', 'Real Sociedad', 'Celta Vigo'}"
$ws.Cells.Item(15, 4).Value = "{'Real Murcia', 'Almeria', 'Real Valladolid'}"
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "2007-2008"
$ws.Cells.Item(16, 3).Value = "{'Real Murcia', 'Real Zaragoza', 'Levante'}"
$ws.Cells.Item(16, 4).Value = "{'Numancia', 'Malaga', 'Sporting Gijon'}"
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "2008-2009"
$ws.Cells.Item(17, 3).Value = "{'Recreativo Huelva', 'Numancia', 'Real Betis'}"
$ws.Cells.Item(17, 4).Value = "{'Xerez', 'Real Zaragoza', 'Tenerife'}"
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "2009-2010"
$ws.Cells.Item(18, 3).Value = "{'Xerez', 'Tenerife', 'Real Valladolid'}"
$ws.Cells.Item(18, 4).Value = "{'Hércules', 'Real Sociedad', 'Levante'}"
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "2010-2011"
$ws.Cells.Item(19, 3).Value = "{'Hércules', 'Almeria', 'Deportivo La C…'}"
$ws.Cells.Item(19, 4).Value = "{'Granada', 'Rayo Vallecano', 'Real Betis'}"
$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = "2011-2012"
$ws.Cells.Item(20, 3).Value = "{'Racing Santander', 'Villarreal', 'Sporting Gijon'}"
$ws.Cells.Item(20, 4).Value = "{'Celta Vigo', 'Real Valladolid', 'Deportivo La C…'}"
$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = "2012-2013"
$ws.Cells.Item(21, 3).Value = "{'Real Zaragoza', 'Mallorca', 'Deportivo La C…'}"
$ws.Cells.Item(21, 4).Value = "{'Villarreal', 'Elche', 'Almeria'}"
$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = "2013-2014"
$ws.Cells.Item(22, 3).Value = "{'Osasuna', 'Real Valladolid', 'Real Betis'}"
$ws.Cells.Item(22, 4).Value = "{'Eibar', 'Córdoba', 'Deportivo La C…'}"
$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = "2014-2015"
$ws.Cells.Item(23, 3).Value = "{'Elche', 'Córdoba', 'Almeria'}"
$ws.Cells.Item(23, 4).Value = "{'Las Palmas', 'Sporting Gijon', 'Real Betis'}"
$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).Value = "2015-2016"
$ws.Cells.Item(24, 3).Value = "{'Getafe', 'Rayo Vallecano', 'Levante'}"
$ws.Cells.Item(24, 4).Value = "{'Alaves', 'Leganes', 'Osasuna'}"
$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).Value = "2016-2017"
$ws.Cells.Item(25, 3).Value = "{'Granada', 'Osasuna', 'Sporting Gijon'}"
$ws.Cells.Item(25, 4).Value = "{'Levante', 'Getafe', 'Girona'}"
$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).Value = "2017-2018"
$ws.Cells.Item(26, 3).Value = "{'Las Palmas', 'Malaga', 'Deportivo La C…'}"
$ws.Cells.Item(26, 4).Value = "{'Huesca', 'Rayo Vallecano', 'Real Valladolid'}"
$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).Value = "2018-2019"
$ws.Cells.Item(27, 3).Value = "{'Huesca', 'Rayo Vallecano', 'Girona'}"
$ws.Cells.Item(27, 4).Value = "{'Granada', 'Osasuna', 'Mallorca'}"
$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = "2019-2020"
$ws.Cells.Item(28, 3).Value = "{'Espanyol', 'Leganes', 'Mallorca'}"
$ws.Cells.Item(28, 4).Value = "{'Elche', 'Cadiz', 'Huesca'}"
$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = "2020-2021"
$ws.Cells.Item(29, 3).Value = "{'Huesca', 'Real Valladolid', 'Eibar'}"
$ws.Cells.Item(29, 4).Value = "{'Espanyol', 'Rayo Vallecano', 'Mallorca'}"
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = "2021-2022"
$ws.Cells.Item(30, 3).Value = "{'Alaves', 'Granada', 'Levante'}"
$ws.Cells.Item(30, 4).Value = "{'Almeria', 'Real Valladolid', 'Girona'}"
$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = "2022-2023"
$ws.Cells.Item(31, 3).Value = "{'Espanyol', 'Elche', 'Real Valladolid'}"
$ws.Cells.Item(31, 4).Value = "{'Las Palmas', 'Granada', 'Alaves'}"

# Style column A (season index) to match the bold/bordered/centered look used elsewhere
$aColRange = $ws.Range("A2:A31")
$aColRange.Font.Bold = $true
$aColRange.Borders.LineStyle = 1
$aColRange.HorizontalAlignment = -4108
$aColRange.VerticalAlignment = -4160

$ws.Range("A1").Select()
